# Auto-applies the cell-value updates from the commit diff
# (profit recalculation sweep across the Anima_Profits workbook sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6925.5
$ws.Range("I62").Value = 2319.75
$ws.Range("K62").Value = 2319.75
$ws.Range("M62").Value = -1695.75
$ws.Range("H65").Value = 6925.5
$ws.Range("I65").Value = 2319.75
$ws.Range("K65").Value = 11598.75
$ws.Range("M65").Value = -8478.75
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H98").Value = 1063.7407
$ws.Range("I98").Value = 866.625
$ws.Range("J98").Value = 2640.6667
$ws.Range("K98").Value = 866.625
$ws.Range("L98").Value = 2640.6667
$ws.Range("M98").Value = 631.375
$ws.Range("N98").Value = -5636.6667
$ws.Range("H122").Value = 1063.7407
$ws.Range("I122").Value = 866.625
$ws.Range("J122").Value = 2640.6667
$ws.Range("K122").Value = 2599.875
$ws.Range("L122").Value = 7922.000100000001
$ws.Range("M122").Value = -149.875
$ws.Range("N122").Value = -12822.0001
$ws.Range("H129").Value = 1397.186
$ws.Range("J129").Value = 1809.1333
$ws.Range("L129").Value = 5427.3999
$ws.Range("N129").Value = -15427.3999
$ws.Range("H132").Value = 3774.4167
$ws.Range("I132").Value = 3771.2856
$ws.Range("J132").Value = 3796.3333
$ws.Range("K132").Value = 11313.8568
$ws.Range("L132").Value = 11388.9999
$ws.Range("M132").Value = -8783.856800000001
$ws.Range("N132").Value = -16448.9999
$ws.Range("H141").Value = 4011.4546
$ws.Range("I141").Value = 2089.7856
$ws.Range("J141").Value = 7374.375
$ws.Range("K141").Value = 6269.3568
$ws.Range("L141").Value = 22123.125
$ws.Range("M141").Value = -1089.3568
$ws.Range("N141").Value = -32483.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5416.23
$ws.Range("I32").Value = 3338.5784
$ws.Range("J32").Value = 15560.059
$ws.Range("K32").Value = 3338.5784
$ws.Range("L32").Value = 15560.059
$ws.Range("M32").Value = -3051.5784
$ws.Range("N32").Value = -16134.059
$ws.Range("H61").Value = 7579580
$ws.Range("I61").Value = 15153650
$ws.Range("J61").Value = 5509.636
$ws.Range("K61").Value = 15153650
$ws.Range("L61").Value = 5509.636
$ws.Range("M61").Value = -15153438
$ws.Range("N61").Value = -5933.636
$ws.Range("H136").Value = 7579580
$ws.Range("I136").Value = 15153650
$ws.Range("J136").Value = 5509.636
$ws.Range("K136").Value = 45460950
$ws.Range("L136").Value = 16528.908
$ws.Range("M136").Value = -45458400
$ws.Range("N136").Value = -21628.908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 625
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2234.8276
$ws.Range("I132").Value = 1926.8889
$ws.Range("J132").Value = 2738.7273
$ws.Range("K132").Value = 5780.6667
$ws.Range("L132").Value = 8216.1819
$ws.Range("M132").Value = -3250.6667
$ws.Range("N132").Value = -13276.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 604.9167
$ws.Range("I113").Value = 555.1667
$ws.Range("J113").Value = 654.6667
$ws.Range("K113").Value = 1665.5001
$ws.Range("L113").Value = 1964.0001
$ws.Range("M113").Value = 504.4999
$ws.Range("N113").Value = -6304.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 23800
$ws.Range("J96").Value = 23800
$ws.Range("L96").Value = 23800
$ws.Range("N96").Value = -29292
$ws.Range("H102").Value = 899.375
$ws.Range("I102").Value = 872.65216
$ws.Range("J102").Value = 1514
$ws.Range("K102").Value = 872.65216
$ws.Range("L102").Value = 1514
$ws.Range("M102").Value = 749.34784
$ws.Range("N102").Value = -4758
$ws.Range("H132").Value = 2818.6428
$ws.Range("I132").Value = 1960.7368
$ws.Range("K132").Value = 5882.2104
$ws.Range("M132").Value = -3352.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5819.091
$ws.Range("I7").Value = 6088
$ws.Range("J7").Value = 5242.857
$ws.Range("K7").Value = 6088
$ws.Range("L7").Value = 5242.857
$ws.Range("M7").Value = -5976
$ws.Range("N7").Value = -5466.857
$ws.Range("H22").Value = 20493.445
$ws.Range("I22").Value = 688
$ws.Range("J22").Value = 45250.25
$ws.Range("K22").Value = 688
$ws.Range("L22").Value = 45250.25
$ws.Range("M22").Value = -393
$ws.Range("N22").Value = -45840.25
$ws.Range("H27").Value = 20493.445
$ws.Range("I27").Value = 688
$ws.Range("J27").Value = 45250.25
$ws.Range("K27").Value = 688
$ws.Range("L27").Value = 45250.25
$ws.Range("M27").Value = -581
$ws.Range("N27").Value = -45464.25
$ws.Range("H40").Value = 1521.7273
$ws.Range("I40").Value = 1340.6842
$ws.Range("J40").Value = 2668.3333
$ws.Range("K40").Value = 1340.6842
$ws.Range("L40").Value = 2668.3333
$ws.Range("M40").Value = -1204.6842
$ws.Range("N40").Value = -2940.3333
$ws.Range("H61").Value = 3165.9473
$ws.Range("I61").Value = 2439.5
$ws.Range("K61").Value = 2439.5
$ws.Range("M61").Value = -2237.5
$ws.Range("H113").Value = 3165.9473
$ws.Range("I113").Value = 2439.5
$ws.Range("K113").Value = 2439.5
$ws.Range("M113").Value = -269.5
$ws.Range("H126").Value = 5819.091
$ws.Range("I126").Value = 6088
$ws.Range("J126").Value = 5242.857
$ws.Range("K126").Value = 18264
$ws.Range("L126").Value = 15728.571
$ws.Range("M126").Value = -15794
$ws.Range("N126").Value = -20668.571
$ws.Range("H132").Value = 3683.5483
$ws.Range("I132").Value = 3652.2666
$ws.Range("J132").Value = 3712.875
$ws.Range("K132").Value = 10956.7998
$ws.Range("L132").Value = 11138.625
$ws.Range("M132").Value = -8426.799800000001
$ws.Range("N132").Value = -16198.625
$ws.Range("H133").Value = 52950.668
$ws.Range("J133").Value = 52950.668
$ws.Range("L133").Value = 52950.668
$ws.Range("N133").Value = -58010.668
$ws.Range("H136").Value = 1934.8
$ws.Range("I136").Value = 1371.1428
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 4113.428400000001
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -1563.428400000001
$ws.Range("N136").Value = -14850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5030.6665
$ws.Range("I81").Value = 7314
$ws.Range("J81").Value = 3399.7144
$ws.Range("K81").Value = 14628
$ws.Range("L81").Value = 6799.4288
$ws.Range("M81").Value = -13567
$ws.Range("N81").Value = -8921.4288
$ws.Range("H84").Value = 5030.6665
$ws.Range("I84").Value = 7314
$ws.Range("J84").Value = 3399.7144
$ws.Range("K84").Value = 73140
$ws.Range("L84").Value = 33997.144
$ws.Range("M84").Value = -67836
$ws.Range("N84").Value = -44605.144
$ws.Range("H132").Value = 7294721.5
$ws.Range("I132").Value = 3476
$ws.Range("J132").Value = 15353467
$ws.Range("K132").Value = 10428
$ws.Range("L132").Value = 46060401
$ws.Range("M132").Value = -7898
$ws.Range("N132").Value = -46065461
